$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.458.76"
$ws.Range("E2").Value = "  -6.70%  "
$ws.Range("D3").Value = "2.587.02"
$ws.Range("E3").Value = "  -0.97%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").Value = "'300.29"
$ws.Range("E5").Value = "  -2.85%  "
$ws.Range("D6").Value = "'96.03"
$ws.Range("E6").Value = "  -4.72%  "
$ws.Range("D7").Value = "'0.576"
$ws.Range("E7").Value = "  -4.80%  "
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("D9").Value = "'0.556"
$ws.Range("E9").Value = "  -4.51%  "
$ws.Range("D10").Value = "'36.80"
$ws.Range("E10").Value = "  -6.95%  "
$ws.Range("D11").Value = "'0.0813"
$ws.Range("E11").Value = "  -4.31%  "
$ws.Range("E12").Value = "  -4.97%  "
$ws.Range("D13").Value = "2.981.20"
$ws.Range("E13").Value = "  -0.74%  "
$ws.Range("E14").Value = "  +0.83%  "
$ws.Range("D15").Value = "2.581.66"
$ws.Range("E15").Value = "  -1.13%  "
$ws.Range("D16").Value = "'0.889"
$ws.Range("E16").Value = "  -4.26%  "
$ws.Range("D17").Value = "'14.31"
$ws.Range("E17").Value = "  -5.20%  "
$ws.Range("D18").Value = "43.417.08"
$ws.Range("E18").Value = "  -7.10%  "
$ws.Range("D19").Value = "'6.65"
$ws.Range("E19").Value = "  -2.42%  "
$ws.Range("E20").Value = "  -4.56%  "
$ws.Range("D21").Value = "'12.30"
$ws.Range("E21").Value = "  -6.48%  "
$ws.Range("D22").Value = "'72.75"
$ws.Range("E22").Value = "  +1.13%  "
$ws.Range("D23").Value = "'265.93"
$ws.Range("E23").Value = "  -4.64%  "
$ws.Range("E24").Value = "  +1.00%  "
$ws.Range("E25").Value = "  -5.00%  "
$ws.Range("D26").Value = "'29.38"
$ws.Range("E26").Value = "  +0.95%  "
$ws.Range("E27").Value = "  -0.30%  "
$ws.Range("D28").Value = "'10.25"
$ws.Range("E28").Value = "  -4.44%  "
$ws.Range("E29").Value = "  -4.60%  "
$ws.Range("D30").Value = "'37.35"
$ws.Range("E30").Value = "  -5.14%  "
$ws.Range("D31").Value = "'6.03"
$ws.Range("E31").Value = "  -5.49%  "
$ws.Range("D32").Value = "'3.58"
$ws.Range("E32").Value = "  -1.72%  "
$ws.Range("D33").Value = "'2.24"
$ws.Range("E33").Value = "  +0.75%  "
$ws.Range("D34").Value = "'152.11"
$ws.Range("E34").Value = "  -0.10%  "
$ws.Range("D35").Value = "'2.78"
$ws.Range("E35").Value = "  -2.16%  "
$ws.Range("D36").Value = "'0.0812"
$ws.Range("E36").Value = "  -3.96%  "
$ws.Range("E37").Value = "  -5.28%  "
$ws.Range("B38").Value = "Stellar"
$ws.Range("C38").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D38").Value = "'0.121"
$ws.Range("E38").Value = "  -1.52%  "
$ws.Range("B39").Value = "EnergySwap"
$ws.Range("C39").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D39").Value = "'24.33"
$ws.Range("E39").Value = "  +4.61%  "
$ws.Range("D40").Value = "'16.62"
$ws.Range("E40").Value = "  +1.88%  "
$ws.Range("D41").Value = "'3.54"
$ws.Range("E41").Value = "  -3.37%  "
$ws.Range("E42").Value = "  -6.31%  "
$ws.Range("E43").Value = "  -6.45%  "
$ws.Range("D44").Value = "2.044.29"
$ws.Range("E44").Value = "  -4.54%  "
$ws.Range("D45").Value = "'0.997"
$ws.Range("E45").Value = "  -0.13%  "
$ws.Range("D46").Value = "'88.54"
$ws.Range("E46").Value = "  -4.94%  "
$ws.Range("D47").Value = "'9.02"
$ws.Range("E47").Value = "  -5.87%  "
$ws.Range("B48").Value = "ApeXProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D48").Value = "'1.61"
$ws.Range("E48").Value = "  +0.52%  "
$ws.Range("B49").Value = "RocketPoolETH"
$ws.Range("C49").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D49").Value = "2.836.62"
$ws.Range("E49").Value = "  -0.81%  "
$ws.Range("D50").Value = "'105.52"
$ws.Range("E50").Value = "  -4.13%  "
$ws.Range("D51").Value = "'0.190"
$ws.Range("E51").Value = "  -5.94%  "
